$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D38").Value = 44895
$ws.Range("K38").Value = 7000
$ws.Range("L38").Value = 7500
$ws.Range("M38").Value = 7250
$ws.Range("O38").Value = 'Región de O''Higgins'
$ws.Range("P38").Value = 145
$ws.Range("D39").Value = 44539
$ws.Range("J39").Value = 270
$ws.Range("K39").Value = 5000
$ws.Range("L39").Value = 5500
$ws.Range("M39").Value = 5222
$ws.Range("N39").Value = '$/caja 60 unidades'
$ws.Range("P39").Value = 87
$ws.Range("Q39").Value = 60
$ws.Range("D40").Value = 44237
$ws.Range("J40").Value = 200
$ws.Range("K40").Value = 8000
$ws.Range("L40").Value = 9000
$ws.Range("M40").Value = 8500
$ws.Range("P40").Value = 142
$ws.Range("D41").Value = 44902
$ws.Range("J41").Value = 430
$ws.Range("K41").Value = 6500
$ws.Range("L41").Value = 7000
$ws.Range("M41").Value = 6733
$ws.Range("N41").Value = '$/caja 50 unidades'
$ws.Range("P41").Value = 135
$ws.Range("Q41").Value = 50
$ws.Range("D42").Value = 44540
$ws.Range("H42").Value = 'Huracán'
$ws.Range("J42").Value = 190
$ws.Range("K42").Value = 6000
$ws.Range("L42").Value = 6500
$ws.Range("M42").Value = 6263
$ws.Range("N42").Value = '$/caja 60 unidades'
$ws.Range("O42").Value = 'Región del Maule'
$ws.Range("P42").Value = 104
$ws.Range("Q42").Value = 60
$ws.Range("D43").Value = 44777
$ws.Range("H43").Value = 'Sin especificar'
$ws.Range("J43").Value = 220
$ws.Range("K43").Value = 17000
$ws.Range("L43").Value = 18000
$ws.Range("M43").Value = 17545
$ws.Range("N43").Value = '$/caja 50 unidades'
$ws.Range("O43").Value = 'Región de Arica y Parinacota'
$ws.Range("P43").Value = 351
$ws.Range("Q43").Value = 50
$ws.Range("D44").Value = 44999
$ws.Range("J44").Value = 120
$ws.Range("K44").Value = 8000
$ws.Range("L44").Value = 8000
$ws.Range("M44").Value = 8000
$ws.Range("O44").Value = 'Región Metropolitana'
$ws.Range("P44").Value = 160
$ws.Range("I45").Value = 'Segunda'
$ws.Range("J45").Value = 100
$ws.Range("K45").Value = 7500
$ws.Range("L45").Value = 7500
$ws.Range("M45").Value = 7500
$ws.Range("P45").Value = 150
$ws.Range("D46").Value = 44831
$ws.Range("I46").Value = 'Primera'
$ws.Range("J46").Value = 270
$ws.Range("K46").Value = 18000
$ws.Range("L46").Value = 19000
$ws.Range("M46").Value = 18444
$ws.Range("O46").Value = 'Región de Arica y Parinacota'
$ws.Range("P46").Value = 369
$ws.Range("D47").Value = 44372
$ws.Range("J47").Value = 100
$ws.Range("K47").Value = 9000
$ws.Range("L47").Value = 10000
$ws.Range("M47").Value = 9500
$ws.Range("P47").Value = 190
$ws.Range("D48").Value = 44370
$ws.Range("K48").Value = 10000
$ws.Range("L48").Value = 11000
$ws.Range("M48").Value = 10500
$ws.Range("P48").Value = 210
$ws.Range("D49").Value = 44533
$ws.Range("J49").Value = 250
$ws.Range("K49").Value = 6000
$ws.Range("L49").Value = 6500
$ws.Range("M49").Value = 6300
$ws.Range("N49").Value = '$/caja 60 unidades'
$ws.Range("P49").Value = 105
$ws.Range("Q49").Value = 60
$ws.Range("D50").Value = 44953
$ws.Range("J50").Value = 50
$ws.Range("K50").Value = 9000
$ws.Range("L50").Value = 9000
$ws.Range("M50").Value = 9000
$ws.Range("N50").Value = '$/caja 50 unidades'
$ws.Range("O50").Value = 'Región de O''Higgins'
$ws.Range("P50").Value = 180
$ws.Range("Q50").Value = 50
$ws.Range("D51").Value = 44894
$ws.Range("J51").Value = 100
$ws.Range("K51").Value = 7000
$ws.Range("L51").Value = 7500
$ws.Range("M51").Value = 7250
$ws.Range("P51").Value = 145
$ws.Range("D52").Value = 44385
$ws.Range("K52").Value = 9000
$ws.Range("L52").Value = 10000
$ws.Range("M52").Value = 9500
$ws.Range("O52").Value = 'Región de Arica y Parinacota'
$ws.Range("P52").Value = 190
$ws.Range("D53").Value = 44596
$ws.Range("J53").Value = 200
$ws.Range("K53").Value = 10000
$ws.Range("L53").Value = 11000
$ws.Range("M53").Value = 10500
$ws.Range("O53").Value = 'Región de O''Higgins'
$ws.Range("P53").Value = 210
$ws.Range("D54").Value = 44875
$ws.Range("J54").Value = 300
$ws.Range("K54").Value = 7500
$ws.Range("L54").Value = 8000
$ws.Range("M54").Value = 7750
$ws.Range("P54").Value = 155
$ws.Range("D55").Value = 44904
$ws.Range("J55").Value = 450
$ws.Range("K55").Value = 6000
$ws.Range("L55").Value = 6500
$ws.Range("M55").Value = 6278
$ws.Range("O55").Value = 'Región Metropolitana'
$ws.Range("P55").Value = 126
$ws.Range("D56").Value = 44203
$ws.Range("J56").Value = 200
$ws.Range("K56").Value = 10000
$ws.Range("L56").Value = 11000
$ws.Range("M56").Value = 10500
$ws.Range("N56").Value = '$/caja 60 unidades'
$ws.Range("O56").Value = 'Región de O''Higgins'
$ws.Range("P56").Value = 175
$ws.Range("Q56").Value = 60
$ws.Range("D57").Value = 44918
$ws.Range("J57").Value = 100
$ws.Range("K57").Value = 7000
$ws.Range("L57").Value = 7500
$ws.Range("M57").Value = 7250
$ws.Range("N57").Value = '$/caja 50 unidades'
$ws.Range("P57").Value = 145
$ws.Range("Q57").Value = 50
$ws.Range("D58").Value = 44580
$ws.Range("K58").Value = 14000
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = 14500
$ws.Range("P58").Value = 290
$ws.Range("D59").Value = 44736
$ws.Range("K59").Value = 16000
$ws.Range("L59").Value = 17000
$ws.Range("M59").Value = 16500
$ws.Range("N59").Value = '$/caja 60 unidades'
$ws.Range("O59").Value = 'Región de Arica y Parinacota'
$ws.Range("P59").Value = 275
$ws.Range("Q59").Value = 60
$ws.Range("D60").Value = 44642
$ws.Range("J60").Value = 220
$ws.Range("K60").Value = 13000
$ws.Range("L60").Value = 14000
$ws.Range("M60").Value = 13455
$ws.Range("P60").Value = 224
$ws.Range("D61").Value = 44966
$ws.Range("J61").Value = 100
$ws.Range("K61").Value = 8500
$ws.Range("L61").Value = 9000
$ws.Range("M61").Value = 8750
$ws.Range("N61").Value = '$/caja 50 unidades'
$ws.Range("O61").Value = 'Región de O''Higgins'
$ws.Range("P61").Value = 175
$ws.Range("Q61").Value = 50
$ws.Range("D62").Value = 44434
$ws.Range("K62").Value = 12000
$ws.Range("L62").Value = 13000
$ws.Range("M62").Value = 12500
$ws.Range("O62").Value = 'Región de Arica y Parinacota'
$ws.Range("P62").Value = 250
$ws.Range("D63").Value = 44162
$ws.Range("J63").Value = 200
$ws.Range("K63").Value = 6000
$ws.Range("L63").Value = 6500
$ws.Range("M63").Value = 6250
$ws.Range("N63").Value = '$/caja 60 unidades'
$ws.Range("O63").Value = 'Región de O''Higgins'
$ws.Range("P63").Value = 104
$ws.Range("Q63").Value = 60
$ws.Range("D64").Value = 44355
$ws.Range("J64").Value = 100
$ws.Range("K64").Value = 9000
$ws.Range("L64").Value = 10000
$ws.Range("M64").Value = 9500
$ws.Range("N64").Value = '$/caja 50 unidades'
$ws.Range("O64").Value = 'Región de Arica y Parinacota'
$ws.Range("P64").Value = 190
$ws.Range("Q64").Value = 50
$ws.Range("D65").Value = 44546
$ws.Range("J65").Value = 250
$ws.Range("K65").Value = 8000
$ws.Range("L65").Value = 9000
$ws.Range("M65").Value = 8400
$ws.Range("N65").Value = '$/caja 60 unidades'
$ws.Range("O65").Value = 'Región Metropolitana'
$ws.Range("P65").Value = 140
$ws.Range("Q65").Value = 60
$ws.Range("D66").Value = 44817
$ws.Range("H66").Value = 'Huracán'
$ws.Range("J66").Value = 220
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 16000
$ws.Range("M66").Value = 15455
$ws.Range("O66").Value = 'Región de Arica y Parinacota'
$ws.Range("P66").Value = 258
$ws.Range("D67").Value = 44649
$ws.Range("H67").Value = 'Sin especificar'
$ws.Range("K67").Value = 12000
$ws.Range("L67").Value = 13000
$ws.Range("M67").Value = 12455
$ws.Range("P67").Value = 208
$ws.Range("D68").Value = 44293
$ws.Range("J68").Value = 100
$ws.Range("K68").Value = 8000
$ws.Range("L68").Value = 9000
$ws.Range("M68").Value = 8500
$ws.Range("O68").Value = 'Región del Maule'
$ws.Range("P68").Value = 142
$ws.Range("D69").Value = 44980
$ws.Range("J69").Value = 150
$ws.Range("K69").Value = 6500
$ws.Range("L69").Value = 7000
$ws.Range("M69").Value = 6833
$ws.Range("N69").Value = '$/caja 50 unidades'
$ws.Range("O69").Value = 'Región de O''Higgins'
$ws.Range("P69").Value = 137
$ws.Range("Q69").Value = 50
$ws.Range("D70").Value = 44341
$ws.Range("J70").Value = 100
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 10000
$ws.Range("M70").Value = 9500
$ws.Range("O70").Value = 'Región de Arica y Parinacota'
$ws.Range("P70").Value = 190
$ws.Range("D71").Value = 44635
$ws.Range("J71").Value = 220
$ws.Range("K71").Value = 12000
$ws.Range("L71").Value = 13000
$ws.Range("M71").Value = 12545
$ws.Range("N71").Value = '$/caja 60 unidades'
$ws.Range("O71").Value = 'Región Metropolitana'
$ws.Range("P71").Value = 209
$ws.Range("Q71").Value = 60
$ws.Range("D72").Value = 44757
$ws.Range("J72").Value = 100
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 10000
$ws.Range("M72").Value = 9500
$ws.Range("N72").Value = '$/caja 50 unidades'
$ws.Range("O72").Value = 'Región de Arica y Parinacota'
$ws.Range("P72").Value = 190
$ws.Range("Q72").Value = 50
$ws.Range("D73").Value = 44327
$ws.Range("D74").Value = 44995
$ws.Range("K74").Value = 8000
$ws.Range("L74").Value = 8500
$ws.Range("M74").Value = 8250
$ws.Range("O74").Value = 'Región de O''Higgins'
$ws.Range("P74").Value = 165
$ws.Range("D75").Value = 44890
$ws.Range("J75").Value = 450
$ws.Range("K75").Value = 6500
$ws.Range("L75").Value = 7000
$ws.Range("M75").Value = 6778
$ws.Range("P75").Value = 136
$ws.Range("D76").Value = 44187
$ws.Range("J76").Value = 200
$ws.Range("K76").Value = 8000
$ws.Range("L76").Value = 9000
$ws.Range("M76").Value = 8500
$ws.Range("N76").Value = '$/caja 60 unidades'
$ws.Range("P76").Value = 142
$ws.Range("Q76").Value = 60
$ws.Range("D77").Value = 44495
$ws.Range("J77").Value = 380
$ws.Range("L77").Value = 8500
$ws.Range("M77").Value = 8237
$ws.Range("N77").Value = '$/caja 50 unidades'
$ws.Range("O77").Value = 'Región de Arica y Parinacota'
$ws.Range("P77").Value = 165
$ws.Range("Q77").Value = 50
$ws.Range("J78").Value = 350
$ws.Range("K78").Value = 10000
$ws.Range("L78").Value = 11000
$ws.Range("M78").Value = 10429
$ws.Range("O78").Value = 'Región de O''Higgins'
$ws.Range("P78").Value = 209
$ws.Range("D79").Value = 44560
$ws.Range("J79").Value = 100
$ws.Range("K79").Value = 6000
$ws.Range("L79").Value = 7000
$ws.Range("M79").Value = 6500
$ws.Range("N79").Value = '$/caja 60 unidades'
$ws.Range("O79").Value = 'Región del Maule'
$ws.Range("P79").Value = 108
$ws.Range("Q79").Value = 60
$ws.Range("D80").Value = 44708
$ws.Range("K80").Value = 18000
$ws.Range("L80").Value = 19000
$ws.Range("M80").Value = 18500
$ws.Range("N80").Value = '$/caja 50 unidades'
$ws.Range("O80").Value = 'Región de Arica y Parinacota'
$ws.Range("P80").Value = 370
$ws.Range("Q80").Value = 50
$ws.Range("D81").Value = 44467
$ws.Range("K81").Value = 14000
$ws.Range("L81").Value = 15000
$ws.Range("M81").Value = 14500
$ws.Range("P81").Value = 290
$ws.Range("D82").Value = 44936
$ws.Range("K82").Value = 7000
$ws.Range("L82").Value = 8000
$ws.Range("M82").Value = 7500
$ws.Range("O82").Value = 'Región de O''Higgins'
$ws.Range("P82").Value = 150
$ws.Range("D83").Value = 44589
$ws.Range("J83").Value = 300
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 11000
$ws.Range("M83").Value = 10500
$ws.Range("N83").Value = '$/caja 60 unidades'
$ws.Range("P83").Value = 175
$ws.Range("Q83").Value = 60
$ws.Range("D84").Value = 44623
$ws.Range("J84").Value = 220
$ws.Range("M84").Value = 10455
$ws.Range("O84").Value = 'Región de Arica y Parinacota'
$ws.Range("P84").Value = 174
$ws.Range("D85").Value = 44670
$ws.Range("J85").Value = 180
$ws.Range("M85").Value = 10444
$ws.Range("N85").Value = '$/caja 50 unidades'
$ws.Range("O85").Value = 'Región Metropolitana'
$ws.Range("P85").Value = 209
$ws.Range("Q85").Value = 50
$ws.Range("D86").Value = 44518
$ws.Range("J86").Value = 450
$ws.Range("K86").Value = 6500
$ws.Range("L86").Value = 7000
$ws.Range("M86").Value = 6722
$ws.Range("N86").Value = '$/caja 60 unidades'
$ws.Range("O86").Value = 'Región de O''Higgins'
$ws.Range("P86").Value = 112
$ws.Range("Q86").Value = 60
$ws.Range("D87").Value = 44944
$ws.Range("J87").Value = 350
$ws.Range("K87").Value = 9000
$ws.Range("L87").Value = 10000
$ws.Range("M87").Value = 9429
$ws.Range("N87").Value = '$/caja 50 unidades'
$ws.Range("O87").Value = 'Región Metropolitana'
$ws.Range("P87").Value = 189
$ws.Range("Q87").Value = 50
$ws.Range("D88").Value = 44883
$ws.Range("J88").Value = 100
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 11000
$ws.Range("M88").Value = 10500
$ws.Range("O88").Value = 'Región de O''Higgins'
$ws.Range("P88").Value = 210
$ws.Range("D89").Value = 44791
$ws.Range("J89").Value = 220
$ws.Range("K89").Value = 21000
$ws.Range("L89").Value = 22000
$ws.Range("M89").Value = 21455
$ws.Range("N89").Value = '$/caja 60 unidades'
$ws.Range("O89").Value = 'Región de Arica y Parinacota'
$ws.Range("P89").Value = 358
$ws.Range("Q89").Value = 60
$ws.Range("D90").Value = 44981
$ws.Range("J90").Value = 100
$ws.Range("K90").Value = 7500
$ws.Range("L90").Value = 8000
$ws.Range("M90").Value = 7750
$ws.Range("N90").Value = '$/caja 50 unidades'
$ws.Range("O90").Value = 'Región de O''Higgins'
$ws.Range("P90").Value = 155
$ws.Range("Q90").Value = 50
$ws.Range("D91").Value = 44323
$ws.Range("K91").Value = 9000
$ws.Range("L91").Value = 10000
$ws.Range("M91").Value = 9500
$ws.Range("O91").Value = 'Región de Arica y Parinacota'
$ws.Range("P91").Value = 190
$ws.Range("D92").Value = 44526
$ws.Range("K92").Value = 7000
$ws.Range("L92").Value = 8000
$ws.Range("M92").Value = 7500
$ws.Range("P92").Value = 150
$ws.Range("D93").Value = 44223
$ws.Range("K93").Value = 9000
$ws.Range("L93").Value = 10000
$ws.Range("M93").Value = 9500
$ws.Range("N93").Value = '$/caja 60 unidades'
$ws.Range("O93").Value = 'Región de O''Higgins'
$ws.Range("P93").Value = 158
$ws.Range("Q93").Value = 60
$ws.Range("D94").Value = 44624
$ws.Range("J94").Value = 150
$ws.Range("K94").Value = 11000
$ws.Range("L94").Value = 12000
$ws.Range("M94").Value = 11467
$ws.Range("O94").Value = 'Región Metropolitana'
$ws.Range("P94").Value = 191
$ws.Range("D95").Value = 44357
$ws.Range("J95").Value = 100
$ws.Range("K95").Value = 8000
$ws.Range("L95").Value = 9000
$ws.Range("M95").Value = 8500
$ws.Range("N95").Value = '$/caja 50 unidades'
$ws.Range("O95").Value = 'Región de Arica y Parinacota'
$ws.Range("P95").Value = 170
$ws.Range("Q95").Value = 50
$ws.Range("D96").Value = 44882
$ws.Range("K96").Value = 9000
$ws.Range("L96").Value = 10000
$ws.Range("M96").Value = 9500
$ws.Range("O96").Value = 'Región de O''Higgins'
$ws.Range("P96").Value = 190
$ws.Range("D97").Value = 44217
$ws.Range("J97").Value = 200
$ws.Range("N97").Value = '$/caja 60 unidades'
$ws.Range("O97").Value = 'Región del Maule'
$ws.Range("P97").Value = 158
$ws.Range("Q97").Value = 60
$ws.Range("D98").Value = 44628
$ws.Range("J98").Value = 220
$ws.Range("K98").Value = 11000
$ws.Range("L98").Value = 12000
$ws.Range("M98").Value = 11545
$ws.Range("O98").Value = 'Región de Arica y Parinacota'
$ws.Range("P98").Value = 192
$ws.Range("D99").Value = 44379
$ws.Range("J99").Value = 100
$ws.Range("K99").Value = 9000
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = 9500
$ws.Range("N99").Value = '$/caja 50 unidades'
$ws.Range("P99").Value = 190
$ws.Range("Q99").Value = 50
$ws.Range("D100").Value = 44755
$ws.Range("K100").Value = 8000
$ws.Range("L100").Value = 8500
$ws.Range("M100").Value = 8250
$ws.Range("P100").Value = 165
$ws.Range("D101").Value = 44484
$ws.Range("J101").Value = 450
$ws.Range("K101").Value = 12000
$ws.Range("L101").Value = 13000
$ws.Range("M101").Value = 12556
$ws.Range("O101").Value = 'Región de O''Higgins'
$ws.Range("P101").Value = 251
$ws.Range("D102").Value = 44778
$ws.Range("J102").Value = 170
$ws.Range("K102").Value = 19000
$ws.Range("L102").Value = 20000
$ws.Range("M102").Value = 19529
$ws.Range("N102").Value = '$/caja 60 unidades'
$ws.Range("O102").Value = 'Región de Arica y Parinacota'
$ws.Range("P102").Value = 325
$ws.Range("Q102").Value = 60
$ws.Range("D103").Value = 44390
$ws.Range("J103").Value = 100
$ws.Range("K103").Value = 9000
$ws.Range("L103").Value = 10000
$ws.Range("M103").Value = 9500
$ws.Range("N103").Value = '$/caja 50 unidades'
$ws.Range("P103").Value = 190
$ws.Range("Q103").Value = 50
$ws.Range("D104").Value = 44950
$ws.Range("O104").Value = 'Región de O''Higgins'
$ws.Range("D105").Value = 44848
$ws.Range("J105").Value = 310
$ws.Range("K105").Value = 17000
$ws.Range("L105").Value = 19000
$ws.Range("M105").Value = 17968
$ws.Range("P105").Value = 359
$ws.Range("D106").Value = 44761
$ws.Range("J106").Value = 220
$ws.Range("K106").Value = 9000
$ws.Range("L106").Value = 10000
$ws.Range("M106").Value = 9545
$ws.Range("N106").Value = '$/caja 60 unidades'
$ws.Range("O106").Value = 'Región de Arica y Parinacota'
$ws.Range("P106").Value = 159
$ws.Range("Q106").Value = 60
$ws.Range("D107").Value = 44678
$ws.Range("J107").Value = 100
$ws.Range("M107").Value = 9500
$ws.Range("N107").Value = '$/caja 50 unidades'
$ws.Range("O107").Value = 'Región Metropolitana'
$ws.Range("P107").Value = 190
$ws.Range("Q107").Value = 50
$ws.Range("D108").Value = 44630
$ws.Range("J108").Value = 200
$ws.Range("K108").Value = 10000
$ws.Range("L108").Value = 11000
$ws.Range("M108").Value = 10500
$ws.Range("N108").Value = '$/caja 60 unidades'
$ws.Range("O108").Value = 'Región del Maule'
$ws.Range("P108").Value = 175
$ws.Range("Q108").Value = 60
$ws.Range("D109").Value = 44602
$ws.Range("J109").Value = 170
$ws.Range("K109").Value = 7000
$ws.Range("L109").Value = 7500
$ws.Range("M109").Value = 7235
$ws.Range("N109").Value = '$/caja 50 unidades'
$ws.Range("P109").Value = 145
$ws.Range("Q109").Value = 50
$ws.Range("D110").Value = 44334
$ws.Range("J110").Value = 100
$ws.Range("K110").Value = 11000
$ws.Range("L110").Value = 12000
$ws.Range("M110").Value = 11500
$ws.Range("O110").Value = 'Región de Arica y Parinacota'
$ws.Range("P110").Value = 230
$ws.Range("D111").Value = 44952
$ws.Range("K111").Value = 9000
$ws.Range("L111").Value = 10000
$ws.Range("M111").Value = 9500
$ws.Range("O111").Value = 'Región de O''Higgins'
$ws.Range("P111").Value = 190
$ws.Range("D112").Value = 44397
$ws.Range("K112").Value = 8000
$ws.Range("L112").Value = 9000
$ws.Range("M112").Value = 8500
$ws.Range("O112").Value = 'Región de Arica y Parinacota'
$ws.Range("P112").Value = 170
$ws.Range("D113").Value = 44505
$ws.Range("J113").Value = 350
$ws.Range("K113").Value = 6500
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = 6714
$ws.Range("N113").Value = '$/caja 60 unidades'
$ws.Range("O113").Value = 'Región del Maule'
$ws.Range("P113").Value = 112
$ws.Range("Q113").Value = 60
$ws.Range("D114").Value = 44659
$ws.Range("K114").Value = 8500
$ws.Range("L114").Value = 9000
$ws.Range("M114").Value = 8714
$ws.Range("N114").Value = '$/caja 50 unidades'
$ws.Range("O114").Value = 'Región Metropolitana'
$ws.Range("P114").Value = 174
$ws.Range("Q114").Value = 50
$ws.Range("D115").Value = 44358
$ws.Range("J115").Value = 100
$ws.Range("K115").Value = 9000
$ws.Range("L115").Value = 10000
$ws.Range("M115").Value = 9500
$ws.Range("O115").Value = 'Región de Arica y Parinacota'
$ws.Range("P115").Value = 190
$ws.Range("D116").Value = 44230
$ws.Range("J116").Value = 150
$ws.Range("M116").Value = 9333
$ws.Range("N116").Value = '$/caja 60 unidades'
$ws.Range("O116").Value = 'Región de O''Higgins'
$ws.Range("P116").Value = 156
$ws.Range("Q116").Value = 60
$ws.Range("D117").Value = 44785
$ws.Range("J117").Value = 100
$ws.Range("K117").Value = 22000
$ws.Range("L117").Value = 23000
$ws.Range("M117").Value = 22500
$ws.Range("N117").Value = '$/caja 50 unidades'
$ws.Range("O117").Value = 'Región de Arica y Parinacota'
$ws.Range("P117").Value = 450
$ws.Range("Q117").Value = 50
$ws.Range("D118").Value = 44638
$ws.Range("J118").Value = 250
$ws.Range("K118").Value = 15000
$ws.Range("L118").Value = 16000
$ws.Range("M118").Value = 15400
$ws.Range("N118").Value = '$/caja 60 unidades'
$ws.Range("O118").Value = 'Región Metropolitana'
$ws.Range("P118").Value = 257
$ws.Range("Q118").Value = 60
$ws.Range("D119").Value = 44243
$ws.Range("J119").Value = 100
$ws.Range("K119").Value = 10000
$ws.Range("L119").Value = 11000
$ws.Range("M119").Value = 10500
$ws.Range("O119").Value = 'Región de O''Higgins'
$ws.Range("P119").Value = 175
$ws.Range("D120").Value = 44460
$ws.Range("N120").Value = '$/caja 50 unidades'
$ws.Range("O120").Value = 'Región de Arica y Parinacota'
$ws.Range("P120").Value = 210
$ws.Range("Q120").Value = 50
$ws.Range("D121").Value = 44775
$ws.Range("J121").Value = 180
$ws.Range("K121").Value = 17000
$ws.Range("L121").Value = 18000
$ws.Range("M121").Value = 17556
$ws.Range("P121").Value = 351
$ws.Range("D122").Value = 44971
$ws.Range("J122").Value = 100
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 8000
$ws.Range("M122").Value = 7750
$ws.Range("O122").Value = 'Región de O''Higgins'
$ws.Range("P122").Value = 155
$ws.Range("D123").Value = 44729
$ws.Range("J123").Value = 180
$ws.Range("K123").Value = 13000
$ws.Range("L123").Value = 14000
$ws.Range("M123").Value = 13556
$ws.Range("O123").Value = 'Región de Arica y Parinacota'
$ws.Range("P123").Value = 271
$ws.Range("D124").Value = 44295
$ws.Range("J124").Value = 100
$ws.Range("K124").Value = 9000
$ws.Range("L124").Value = 10000
$ws.Range("M124").Value = 9500
$ws.Range("N124").Value = '$/caja 60 unidades'
$ws.Range("O124").Value = 'Región de O''Higgins'
$ws.Range("P124").Value = 158
$ws.Range("Q124").Value = 60
$ws.Range("D125").Value = 44299
$ws.Range("K125").Value = 7000
$ws.Range("L125").Value = 8000
$ws.Range("M125").Value = 7500
$ws.Range("N125").Value = '$/caja 50 unidades'
$ws.Range("O125").Value = 'Región Metropolitana'
$ws.Range("P125").Value = 150
$ws.Range("Q125").Value = 50
$ws.Range("D126").Value = 44558
$ws.Range("J126").Value = 250
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 6500
$ws.Range("M126").Value = 6200
$ws.Range("N126").Value = '$/caja 60 unidades'
$ws.Range("O126").Value = 'Región de Arica y Parinacota'
$ws.Range("P126").Value = 103
$ws.Range("Q126").Value = 60
$ws.Range("D127").Value = 44265
$ws.Range("J127").Value = 200
$ws.Range("K127").Value = 7000
$ws.Range("L127").Value = 8000
$ws.Range("M127").Value = 7500
$ws.Range("O127").Value = 'Región de O''Higgins'
$ws.Range("P127").Value = 125
$ws.Range("D128").Value = 44880
$ws.Range("J128").Value = 100
$ws.Range("K128").Value = 10000
$ws.Range("L128").Value = 11000
$ws.Range("M128").Value = 10500
$ws.Range("N128").Value = '$/caja 50 unidades'
$ws.Range("P128").Value = 210
$ws.Range("Q128").Value = 50
$ws.Range("D129").Value = 44645
$ws.Range("J129").Value = 220
$ws.Range("K129").Value = 11000
$ws.Range("L129").Value = 12000
$ws.Range("M129").Value = 11545
$ws.Range("N129").Value = '$/caja 60 unidades'
$ws.Range("O129").Value = 'Región Metropolitana'
$ws.Range("P129").Value = 192
$ws.Range("Q129").Value = 60
$ws.Range("D130").Value = 44932
$ws.Range("J130").Value = 450
$ws.Range("K130").Value = 10000
$ws.Range("L130").Value = 11000
$ws.Range("M130").Value = 10556
$ws.Range("O130").Value = 'Región de Arica y Parinacota'
$ws.Range("P130").Value = 176
$ws.Range("D131").Value = 44616
$ws.Range("J131").Value = 150
$ws.Range("K131").Value = 9000
$ws.Range("L131").Value = 10000
$ws.Range("M131").Value = 9333
$ws.Range("N131").Value = '$/caja 50 unidades'
$ws.Range("O131").Value = 'Región de O''Higgins'
$ws.Range("P131").Value = 187
$ws.Range("Q131").Value = 50
$ws.Range("D132").Value = 45001
$ws.Range("K132").Value = 8000
$ws.Range("L132").Value = 8500
$ws.Range("M132").Value = 8233
$ws.Range("N132").Value = '$/caja 60 unidades'
$ws.Range("O132").Value = 'Región Metropolitana'
$ws.Range("P132").Value = 137
$ws.Range("Q132").Value = 60
